# Apply 2019-04-06 am activity log updates to the "2019" sheet / Table2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$lo = $ws.ListObjects.Item(1)

# --- Pre-format the new date cells (A/B) by copying the existing date-format
#     style from established rows, so the new cells reuse the workbook's existing
#     number-format style instead of minting a duplicate one. ---
$ws.Range("A171").Copy()
$ws.Range("A172:A190").PasteSpecial(-4122)
$ws.Range("B166").Copy()
$ws.Range("B172").PasteSpecial(-4122)
$ws.Range("B180").PasteSpecial(-4122)
$ws.Range("B184").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Write the new rows data ---
# Row 172
$ws.Cells.Item(172, 1).Value = 43558.934108796297
$ws.Cells.Item(172, 2).Value = 43559.220138888886
$ws.Cells.Item(172, 3).Value = "Sleep"
$ws.Cells.Item(172, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 173
$ws.Cells.Item(173, 1).Value = 43559.260416666664
$ws.Cells.Item(173, 3).Value = "Food"
$ws.Cells.Item(173, 4).Value = "Latte"
$ws.Cells.Item(173, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 174
$ws.Cells.Item(174, 1).Value = 43558.78125
$ws.Cells.Item(174, 3).Value = "Food"
$ws.Cells.Item(174, 4).Value = "Enchilada"
$ws.Cells.Item(174, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 175
$ws.Cells.Item(175, 1).Value = 43559.34097222222
$ws.Cells.Item(175, 3).Value = "Food"
$ws.Cells.Item(175, 4).Value = "Whole Wheat Bread"
$ws.Cells.Item(175, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 176
$ws.Cells.Item(176, 1).Value = 43559.53125
$ws.Cells.Item(176, 3).Value = "Food"
$ws.Cells.Item(176, 4).Value = "Enchilada"
$ws.Cells.Item(176, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 177
$ws.Cells.Item(177, 1).Value = 43559.635416666664
$ws.Cells.Item(177, 3).Value = "Food"
$ws.Cells.Item(177, 4).Value = "Kombucha"
$ws.Cells.Item(177, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 178
$ws.Cells.Item(178, 1).Value = 43559.552083333336
$ws.Cells.Item(178, 3).Value = "Food"
$ws.Cells.Item(178, 4).Value = "Banana + nuts"
$ws.Cells.Item(178, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 179
$ws.Cells.Item(179, 1).Value = 43559.423611111109
$ws.Cells.Item(179, 3).Value = "Food"
$ws.Cells.Item(179, 4).Value = "Blueberry Kefir"
$ws.Cells.Item(179, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 180
$ws.Cells.Item(180, 1).Value = 43559.946979166663
$ws.Cells.Item(180, 2).Value = 43560.223611111112
$ws.Cells.Item(180, 3).Value = "Sleep"
$ws.Cells.Item(180, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 181
$ws.Cells.Item(181, 1).Value = 43560.25
$ws.Cells.Item(181, 3).Value = "Food"
$ws.Cells.Item(181, 4).Value = "Latte"
$ws.Cells.Item(181, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 182
$ws.Cells.Item(182, 1).Value = 43560.34375
$ws.Cells.Item(182, 3).Value = "Food"
$ws.Cells.Item(182, 4).Value = "eggs + beans"
$ws.Cells.Item(182, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 183
$ws.Cells.Item(183, 1).Value = 43559.770833333336
$ws.Cells.Item(183, 3).Value = "Food"
$ws.Cells.Item(183, 4).Value = "Pork + banana bread"
$ws.Cells.Item(183, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 184
$ws.Cells.Item(184, 1).Value = 43560.911620370367
$ws.Cells.Item(184, 2).Value = 43561.21875
$ws.Cells.Item(184, 3).Value = "Sleep"
$ws.Cells.Item(184, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 185
$ws.Cells.Item(185, 1).Value = 43560.8125
$ws.Cells.Item(185, 3).Value = "Food"
$ws.Cells.Item(185, 4).Value = "Vegetable Udon"
$ws.Cells.Item(185, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 186
$ws.Cells.Item(186, 1).Value = 43560.520833333336
$ws.Cells.Item(186, 3).Value = "Food"
$ws.Cells.Item(186, 4).Value = "Pesto Pasta"
$ws.Cells.Item(186, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 187
$ws.Cells.Item(187, 1).Value = 43560.604166666664
$ws.Cells.Item(187, 3).Value = "Food"
$ws.Cells.Item(187, 4).Value = "Banana bread"
$ws.Cells.Item(187, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 188
$ws.Cells.Item(188, 1).Value = 43560.6875
$ws.Cells.Item(188, 3).Value = "Food"
$ws.Cells.Item(188, 4).Value = "Banana bread"
$ws.Cells.Item(188, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 189
$ws.Cells.Item(189, 1).Value = 43561.25
$ws.Cells.Item(189, 3).Value = "Food"
$ws.Cells.Item(189, 4).Value = "Latte"
$ws.Cells.Item(189, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# Row 190
$ws.Cells.Item(190, 1).Value = 43561.354166666664
$ws.Cells.Item(190, 3).Value = "Food"
$ws.Cells.Item(190, 4).Value = "eggs"
$ws.Cells.Item(190, 5).Formula = "=IF(Table2[[#This Row],[Activity]]=""Sleep"",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,""NA"")"

# --- Resize the table (ListObject) to cover the newly added rows ---
$lo.Resize($ws.Range("A1:E190"))

# --- Update the sheet view: selection follows the new last row ---
$ws.Range("A191").Select()

